# "New setup added of screenshot, testCaseFolderCreation"
# Two cells that previously held "No" (shared string) now hold "Yes",
# the appointment-fee-adjacent "days" cell in row 4 changes 30 -> 5,
# and the sheet's scroll/selection moves from E13 to a view scrolled to
# column C with J9 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B4: "No" -> "Yes"
$ws.Range("B4").Value = "Yes"

# I4: "No" -> "Yes"
$ws.Range("I4").Value = "Yes"

# L4: 30 -> 5
$ws.Range("L4").Value = 5

# Update the visible/active window: scroll so column C is left-most visible,
# then move the selection to J9 (mirrors the sheetView/selection change).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("J9").Select()
